# Settings workbook update:
#   - Rename header "Tokens Per Hour" -> "Tokens per Hour"
#   - Add a new "Gold per Value" column (G) used to compute gold distribution
#     per token, with a value of 1000 for each existing row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix the casing of the existing "Tokens Per Hour" header.
$ws.Range("E1").Value = "Tokens per Hour"

# New "Gold per Value" setting column, mirroring the other numeric/header cells.
$ws.Range("G1").Value = "Gold per Value"
$ws.Range("G2").Value = 1000.0
$ws.Range("G3").Value = 1000.0

# Match the formatting of the neighboring "Total Time" column (F) so the new
# column looks consistent with the rest of the header/data rows.
$ws.Range("F1:F3").Copy()
$ws.Range("G1:G3").PasteSpecial(-4122)
